$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 112
$ws.Range("H112").Value = 1040994
$ws.Range("I112").Value = 921.4286
$ws.Range("J112").Value = 1445466.6
$ws.Range("K112").Value = 2764.2858
$ws.Range("L112").Value = 4336399.800000001
$ws.Range("M112").Value = -1656.2858
$ws.Range("N112").Value = -4338615.800000001

# Row 132
$ws.Range("H132").Value = 2380.5676
$ws.Range("I132").Value = 1930.8518
$ws.Range("J132").Value = 3594.8
$ws.Range("K132").Value = 5792.555399999999
$ws.Range("L132").Value = 10784.4
$ws.Range("M132").Value = -3262.555399999999
$ws.Range("N132").Value = -15844.4

# Row 137
$ws.Range("H137").Value = 3244.8225
$ws.Range("I137").Value = 1258.5264
$ws.Range("J137").Value = 6389.7915
$ws.Range("K137").Value = 3775.5792
$ws.Range("L137").Value = 19169.3745
$ws.Range("M137").Value = -1225.5792
$ws.Range("N137").Value = -24269.3745

# Row 138
$ws.Range("H138").Value = 2485.1091
$ws.Range("I138").Value = 1261.0646
$ws.Range("J138").Value = 4066.1667
$ws.Range("K138").Value = 3783.1938
$ws.Range("L138").Value = 12198.5001
$ws.Range("M138").Value = 1356.8062
$ws.Range("N138").Value = -22478.5001

# Row 141
$ws.Range("H141").Value = 2923.3572
$ws.Range("I141").Value = 1611.723
$ws.Range("J141").Value = 7410.5264
$ws.Range("K141").Value = 4835.169
$ws.Range("L141").Value = 22231.5792
$ws.Range("M141").Value = 344.8310000000001
$ws.Range("N141").Value = -32591.5792

$ws = $wb.Worksheets.Item("ARM")
# Row 2
$ws.Range("H2").Value = 42882.207
$ws.Range("I2").Value = 56428.61
$ws.Range("J2").Value = 2243
$ws.Range("K2").Value = 56428.61
$ws.Range("L2").Value = 2243
$ws.Range("M2").Value = -56315.61
$ws.Range("N2").Value = -2469

# Row 32
$ws.Range("H32").Value = 10531245
$ws.Range("I32").Value = 11768703
$ws.Range("J32").Value = 12855.6
$ws.Range("K32").Value = 11768703
$ws.Range("L32").Value = 12855.6
$ws.Range("M32").Value = -11768416
$ws.Range("N32").Value = -13429.6

# Row 45
$ws.Range("H45").Value = 1491.3334
$ws.Range("I45").Value = 1373.4286
$ws.Range("J45").Value = 1727.1428
$ws.Range("K45").Value = 1373.4286
$ws.Range("L45").Value = 1727.1428
$ws.Range("M45").Value = -996.4286
$ws.Range("N45").Value = -2481.1428

# Row 74
$ws.Range("H74").Value = 2650.0637
$ws.Range("I74").Value = 490.5122
$ws.Range("J74").Value = 17407
$ws.Range("K74").Value = 490.5122
$ws.Range("L74").Value = 17407
$ws.Range("M74").Value = 383.4878
$ws.Range("N74").Value = -19155

# Row 77
$ws.Range("H77").Value = 2650.0637
$ws.Range("I77").Value = 490.5122
$ws.Range("J77").Value = 17407
$ws.Range("K77").Value = 2452.561
$ws.Range("L77").Value = 87035
$ws.Range("M77").Value = 1915.439
$ws.Range("N77").Value = -95771

# Row 116
$ws.Range("H116").Value = 42882.207
$ws.Range("I116").Value = 56428.61
$ws.Range("J116").Value = 2243
$ws.Range("K116").Value = 56428.61
$ws.Range("L116").Value = 2243
$ws.Range("M116").Value = -54134.61
$ws.Range("N116").Value = -6831

# Row 134
$ws.Range("H134").Value = 49444
$ws.Range("J134").Value = 49444
$ws.Range("L134").Value = 49444
$ws.Range("N134").Value = -59584

$ws = $wb.Worksheets.Item("BSM")
# Row 3
$ws.Range("H3").Value = 42882.207
$ws.Range("I3").Value = 56428.61
$ws.Range("J3").Value = 2243
$ws.Range("K3").Value = 56428.61
$ws.Range("L3").Value = 2243
$ws.Range("M3").Value = -56314.61
$ws.Range("N3").Value = -2471

# Row 26
$ws.Range("H26").Value = 23754.75
$ws.Range("I26").Value = 9000
$ws.Range("J26").Value = 28673
$ws.Range("K26").Value = 9000
$ws.Range("L26").Value = 28673
$ws.Range("M26").Value = -8708
$ws.Range("N26").Value = -29257

# Row 86
$ws.Range("H86").Value = 1996.7354
$ws.Range("I86").Value = 1788.8125
$ws.Range("J86").Value = 2181.5557
$ws.Range("K86").Value = 1788.8125
$ws.Range("L86").Value = 2181.5557
$ws.Range("M86").Value = -665.8125
$ws.Range("N86").Value = -4427.5557

# Row 89
$ws.Range("H89").Value = 1996.7354
$ws.Range("I89").Value = 1788.8125
$ws.Range("J89").Value = 2181.5557
$ws.Range("K89").Value = 8944.0625
$ws.Range("L89").Value = 10907.7785
$ws.Range("M89").Value = -3328.0625
$ws.Range("N89").Value = -22139.7785

# Row 134
$ws.Range("H134").Value = 1075.5
$ws.Range("I134").Value = 759.6539
$ws.Range("J134").Value = 1759.8334
$ws.Range("K134").Value = 2278.9617
$ws.Range("L134").Value = 5279.5002
$ws.Range("M134").Value = 256.0383000000002
$ws.Range("N134").Value = -10349.5002

$ws = $wb.Worksheets.Item("CRP")
# Row 31
$ws.Range("H31").Value = 91636
$ws.Range("I31").Value = 0
$ws.Range("J31").Value = 91636
$ws.Range("K31").Value = 0
$ws.Range("L31").Value = 91636
$ws.Range("M31").ClearContents()
$ws.Range("N31").Value = -92226

# Row 34
$ws.Range("H34").Value = 91636
$ws.Range("I34").Value = 0
$ws.Range("J34").Value = 91636
$ws.Range("K34").Value = 0
$ws.Range("L34").Value = 91636
$ws.Range("M34").ClearContents()
$ws.Range("N34").Value = -92040

# Row 58
$ws.Range("H58").Value = 820.9773
$ws.Range("I58").Value = 554.16394
$ws.Range("J58").Value = 1423.7778
$ws.Range("K58").Value = 554.16394
$ws.Range("L58").Value = 1423.7778
$ws.Range("M58").Value = -351.16394
$ws.Range("N58").Value = -1829.7778

# Row 59
$ws.Range("H59").Value = 15498.75
$ws.Range("I59").Value = 2500
$ws.Range("J59").Value = 19831.666
$ws.Range("K59").Value = 2500
$ws.Range("L59").Value = 19831.666
$ws.Range("M59").Value = -1355
$ws.Range("N59").Value = -22121.666

# Row 99
$ws.Range("H99").Value = 2999.2144
$ws.Range("I99").Value = 2841.1667
$ws.Range("J99").Value = 3117.75
$ws.Range("K99").Value = 2841.1667
$ws.Range("L99").Value = 3117.75
$ws.Range("M99").Value = -1343.1667
$ws.Range("N99").Value = -6113.75

# Row 126
$ws.Range("H126").Value = 2999.2144
$ws.Range("I126").Value = 2841.1667
$ws.Range("J126").Value = 3117.75
$ws.Range("K126").Value = 8523.500100000001
$ws.Range("L126").Value = 9353.25
$ws.Range("M126").Value = -6053.500100000001
$ws.Range("N126").Value = -14293.25

# Row 134
$ws.Range("H134").Value = 1186.9589
$ws.Range("I134").Value = 1132.381
$ws.Range("J134").Value = 1530.8
$ws.Range("K134").Value = 3397.143
$ws.Range("L134").Value = 4592.4
$ws.Range("M134").Value = -862.143
$ws.Range("N134").Value = -9662.4

# Row 136
$ws.Range("H136").Value = 820.9773
$ws.Range("I136").Value = 554.16394
$ws.Range("J136").Value = 1423.7778
$ws.Range("K136").Value = 1662.49182
$ws.Range("L136").Value = 4271.3334
$ws.Range("M136").Value = 887.5081799999998
$ws.Range("N136").Value = -9371.3334

$ws = $wb.Worksheets.Item("CUL")
# Row 74
$ws.Range("H74").Value = 1000
$ws.Range("I74").Value = 1000
$ws.Range("K74").Value = 3000
$ws.Range("M74").Value = -1939

# Row 77
$ws.Range("H77").Value = 1000
$ws.Range("I77").Value = 1000
$ws.Range("K77").Value = 9000
$ws.Range("M77").Value = -3696

# Row 131
$ws.Range("H131").Value = 720.5599999999999
$ws.Range("J131").Value = 965
$ws.Range("L131").Value = 2895
$ws.Range("N131").Value = -12975

$ws = $wb.Worksheets.Item("GSM")
# Row 70
$ws.Range("H70").Value = 3448
$ws.Range("I70").Value = 3381.75
$ws.Range("J70").Value = 3660
$ws.Range("K70").Value = 3381.75
$ws.Range("L70").Value = 3660
$ws.Range("M70").Value = -3111.75
$ws.Range("N70").Value = -4200

# Row 73
$ws.Range("H73").Value = 3448
$ws.Range("I73").Value = 3381.75
$ws.Range("J73").Value = 3660
$ws.Range("K73").Value = 3381.75
$ws.Range("L73").Value = 3660
$ws.Range("M73").Value = -2445.75
$ws.Range("N73").Value = -5532

$ws = $wb.Worksheets.Item("LTW")
# Row 120
$ws.Range("H120").Value = 48666.668
$ws.Range("J120").Value = 48666.668
$ws.Range("L120").Value = 48666.668
$ws.Range("N120").Value = -58342.668

$ws = $wb.Worksheets.Item("WVR")
# Row 132
$ws.Range("H132").Value = 11114936
$ws.Range("I132").Value = 20839048
$ws.Range("J132").Value = 1664.9524
$ws.Range("K132").Value = 62517144
$ws.Range("L132").Value = 4994.857199999999
$ws.Range("M132").Value = -62514614
$ws.Range("N132").Value = -10054.8572
